# Auto-generated edit script applying the Mateus_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 666.3333
$ws.Cells.Item(18, 10).Value = 499
$ws.Cells.Item(18, 12).Value = 499
$ws.Cells.Item(18, 14).Value = -1067
$ws.Cells.Item(80, 8).Value = 646.2857
$ws.Cells.Item(80, 9).Value = 582.3333
$ws.Cells.Item(80, 10).Value = 731.55554
$ws.Cells.Item(80, 11).Value = 1746.9999
$ws.Cells.Item(80, 12).Value = 2194.66662
$ws.Cells.Item(80, 13).Value = -748.9999
$ws.Cells.Item(80, 14).Value = -4190.66662
$ws.Cells.Item(83, 8).Value = 646.2857
$ws.Cells.Item(83, 9).Value = 582.3333
$ws.Cells.Item(83, 10).Value = 731.55554
$ws.Cells.Item(83, 11).Value = 5240.9997
$ws.Cells.Item(83, 12).Value = 6583.99986
$ws.Cells.Item(83, 13).Value = -248.9997000000003
$ws.Cells.Item(83, 14).Value = -16567.99986
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).ClearContents()
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 2019.9656
$ws.Cells.Item(137, 9).Value = 1783.5714
$ws.Cells.Item(137, 11).Value = 5350.7142
$ws.Cells.Item(137, 13).Value = -2800.7142
$ws.Cells.Item(141, 8).Value = 1936.0741
$ws.Cells.Item(141, 9).Value = 1686.8
$ws.Cells.Item(141, 10).Value = 5052
$ws.Cells.Item(141, 11).Value = 5060.4
$ws.Cells.Item(141, 12).Value = 15156
$ws.Cells.Item(141, 13).Value = 119.6000000000004
$ws.Cells.Item(141, 14).Value = -25516

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 58797.223
$ws.Cells.Item(45, 9).Value = 71420.69
$ws.Cells.Item(45, 11).Value = 71420.69
$ws.Cells.Item(45, 13).Value = -71043.69
$ws.Cells.Item(61, 8).Value = 11121951
$ws.Cells.Item(61, 9).Value = 14713523
$ws.Cells.Item(61, 10).Value = 20726.727
$ws.Cells.Item(61, 11).Value = 14713523
$ws.Cells.Item(61, 12).Value = 20726.727
$ws.Cells.Item(61, 13).Value = -14713311
$ws.Cells.Item(61, 14).Value = -21150.727
$ws.Cells.Item(74, 8).Value = 5161.129
$ws.Cells.Item(74, 9).Value = 4302.4165
$ws.Cells.Item(74, 11).Value = 4302.4165
$ws.Cells.Item(74, 13).Value = -3428.4165
$ws.Cells.Item(77, 8).Value = 5161.129
$ws.Cells.Item(77, 9).Value = 4302.4165
$ws.Cells.Item(77, 11).Value = 21512.0825
$ws.Cells.Item(77, 13).Value = -17144.0825
$ws.Cells.Item(97, 8).Value = 1050.3928
$ws.Cells.Item(97, 9).Value = 981
$ws.Cells.Item(97, 11).Value = 981
$ws.Cells.Item(97, 13).Value = -485
$ws.Cells.Item(121, 8).Value = 50000
$ws.Cells.Item(121, 10).Value = 50000
$ws.Cells.Item(121, 12).Value = 50000
$ws.Cells.Item(121, 14).Value = -53494
$ws.Cells.Item(132, 8).Value = 3168.3542
$ws.Cells.Item(132, 9).Value = 3291.2666
$ws.Cells.Item(132, 11).Value = 9873.799800000001
$ws.Cells.Item(132, 13).Value = -7343.799800000001
$ws.Cells.Item(134, 8).Value = 146855.86
$ws.Cells.Item(134, 10).Value = 146855.86
$ws.Cells.Item(134, 12).Value = 146855.86
$ws.Cells.Item(134, 14).Value = -156995.86
$ws.Cells.Item(135, 8).Value = 80319
$ws.Cells.Item(135, 10).Value = 80319
$ws.Cells.Item(135, 12).Value = 80319
$ws.Cells.Item(135, 14).Value = -90459
$ws.Cells.Item(136, 8).Value = 11121951
$ws.Cells.Item(136, 9).Value = 14713523
$ws.Cells.Item(136, 10).Value = 20726.727
$ws.Cells.Item(136, 11).Value = 44140569
$ws.Cells.Item(136, 12).Value = 62180.181
$ws.Cells.Item(136, 13).Value = -44138019
$ws.Cells.Item(136, 14).Value = -67280.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1962.3572
$ws.Cells.Item(105, 9).Value = 1392.8572
$ws.Cells.Item(105, 11).Value = 1392.8572
$ws.Cells.Item(105, 13).Value = 354.1428000000001
$ws.Cells.Item(107, 8).Value = 2112.25
$ws.Cells.Item(107, 9).Value = 2122.4546
$ws.Cells.Item(107, 11).Value = 2122.4546
$ws.Cells.Item(107, 13).Value = -202.4546
$ws.Cells.Item(134, 8).Value = 4402.9487
$ws.Cells.Item(134, 9).Value = 4546.3516
$ws.Cells.Item(134, 11).Value = 13639.0548
$ws.Cells.Item(134, 13).Value = -11104.0548

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4766
$ws.Cells.Item(31, 9).Value = 2492
$ws.Cells.Item(31, 11).Value = 2492
$ws.Cells.Item(31, 13).Value = -2197
$ws.Cells.Item(34, 8).Value = 4766
$ws.Cells.Item(34, 9).Value = 2492
$ws.Cells.Item(34, 11).Value = 2492
$ws.Cells.Item(34, 13).Value = -2290
$ws.Cells.Item(58, 8).Value = 5337.1387
$ws.Cells.Item(58, 9).Value = 3567.8333
$ws.Cells.Item(58, 10).Value = 8875.75
$ws.Cells.Item(58, 11).Value = 3567.8333
$ws.Cells.Item(58, 12).Value = 8875.75
$ws.Cells.Item(58, 13).Value = -3364.8333
$ws.Cells.Item(58, 14).Value = -9281.75
$ws.Cells.Item(74, 8).Value = 42372.5
$ws.Cells.Item(74, 10).Value = 42372.5
$ws.Cells.Item(74, 12).Value = 42372.5
$ws.Cells.Item(74, 14).Value = -44120.5
$ws.Cells.Item(77, 8).Value = 42372.5
$ws.Cells.Item(77, 10).Value = 42372.5
$ws.Cells.Item(77, 12).Value = 127117.5
$ws.Cells.Item(77, 14).Value = -135853.5
$ws.Cells.Item(99, 8).Value = 5051.25
$ws.Cells.Item(99, 9).Value = 5688.9
$ws.Cells.Item(99, 11).Value = 5688.9
$ws.Cells.Item(99, 13).Value = -4190.9
$ws.Cells.Item(107, 8).Value = 2233.875
$ws.Cells.Item(107, 9).Value = 2874.2
$ws.Cells.Item(107, 11).Value = 2874.2
$ws.Cells.Item(107, 13).Value = -954.1999999999998
$ws.Cells.Item(123, 8).Value = 49397.8
$ws.Cells.Item(123, 10).Value = 49397.8
$ws.Cells.Item(123, 12).Value = 49397.8
$ws.Cells.Item(123, 14).Value = -59197.8
$ws.Cells.Item(126, 8).Value = 5051.25
$ws.Cells.Item(126, 9).Value = 5688.9
$ws.Cells.Item(126, 11).Value = 17066.7
$ws.Cells.Item(126, 13).Value = -14596.7
$ws.Cells.Item(133, 8).Value = 72713
$ws.Cells.Item(133, 10).Value = 72713
$ws.Cells.Item(133, 12).Value = 72713
$ws.Cells.Item(133, 14).Value = -77773
$ws.Cells.Item(136, 8).Value = 5337.1387
$ws.Cells.Item(136, 9).Value = 3567.8333
$ws.Cells.Item(136, 10).Value = 8875.75
$ws.Cells.Item(136, 11).Value = 10703.4999
$ws.Cells.Item(136, 12).Value = 26627.25
$ws.Cells.Item(136, 13).Value = -8153.499899999999
$ws.Cells.Item(136, 14).Value = -31727.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 14000880
$ws.Cells.Item(11, 9).Value = 21000126
$ws.Cells.Item(11, 11).Value = 63000378
$ws.Cells.Item(11, 13).Value = -63000238
$ws.Cells.Item(131, 8).Value = 16130593
$ws.Cells.Item(131, 9).Value = 50000844
$ws.Cells.Item(131, 10).Value = 1902.9524
$ws.Cells.Item(131, 11).Value = 150002532
$ws.Cells.Item(131, 12).Value = 5708.857199999999
$ws.Cells.Item(131, 13).Value = -149997492
$ws.Cells.Item(131, 14).Value = -15788.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5556.75
$ws.Cells.Item(80, 9).Value = 5249
$ws.Cells.Item(80, 10).Value = 5659.3335
$ws.Cells.Item(80, 11).Value = 5249
$ws.Cells.Item(80, 12).Value = 5659.3335
$ws.Cells.Item(80, 13).Value = -4251
$ws.Cells.Item(80, 14).Value = -7655.3335
$ws.Cells.Item(83, 8).Value = 5556.75
$ws.Cells.Item(83, 9).Value = 5249
$ws.Cells.Item(83, 10).Value = 5659.3335
$ws.Cells.Item(83, 11).Value = 26245
$ws.Cells.Item(83, 12).Value = 28296.6675
$ws.Cells.Item(83, 13).Value = -21253
$ws.Cells.Item(83, 14).Value = -38280.6675
$ws.Cells.Item(111, 8).Value = 40000
$ws.Cells.Item(111, 10).Value = 40000
$ws.Cells.Item(111, 12).Value = 40000
$ws.Cells.Item(111, 14).Value = -46134
$ws.Cells.Item(122, 8).Value = 6722
$ws.Cells.Item(122, 9).Value = 5352.353
$ws.Cells.Item(122, 11).Value = 16057.059
$ws.Cells.Item(122, 13).Value = -13607.059
$ws.Cells.Item(126, 8).Value = 3075.4546
$ws.Cells.Item(126, 9).Value = 1934.7142
$ws.Cells.Item(126, 11).Value = 5804.142599999999
$ws.Cells.Item(126, 13).Value = -3334.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3022
$ws.Cells.Item(7, 9).Value = 3329.3333
$ws.Cells.Item(7, 11).Value = 3329.3333
$ws.Cells.Item(7, 13).Value = -3217.3333
$ws.Cells.Item(40, 8).Value = 3545.4666
$ws.Cells.Item(40, 9).Value = 2818.2
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 11).Value = 2818.2
$ws.Cells.Item(40, 12).Value = 5000
$ws.Cells.Item(40, 13).Value = -2682.2
$ws.Cells.Item(40, 14).Value = -5272
$ws.Cells.Item(46, 8).Value = 2149.25
$ws.Cells.Item(46, 10).Value = 1949.5
$ws.Cells.Item(46, 12).Value = 1949.5
$ws.Cells.Item(46, 14).Value = -2325.5
$ws.Cells.Item(55, 8).Value = 524.36365
$ws.Cells.Item(55, 9).Value = 778.7143
$ws.Cells.Item(55, 11).Value = 778.7143
$ws.Cells.Item(55, 13).Value = -605.7143
$ws.Cells.Item(61, 8).Value = 64550.75
$ws.Cells.Item(61, 9).Value = 64550.75
$ws.Cells.Item(61, 11).Value = 64550.75
$ws.Cells.Item(61, 13).Value = -64348.75
$ws.Cells.Item(82, 8).Value = 972.6087
$ws.Cells.Item(82, 9).Value = 730.2727
$ws.Cells.Item(82, 10).Value = 1194.75
$ws.Cells.Item(82, 11).Value = 730.2727
$ws.Cells.Item(82, 12).Value = 1194.75
$ws.Cells.Item(82, 13).Value = -369.2727
$ws.Cells.Item(82, 14).Value = -1916.75
$ws.Cells.Item(85, 8).Value = 972.6087
$ws.Cells.Item(85, 9).Value = 730.2727
$ws.Cells.Item(85, 10).Value = 1194.75
$ws.Cells.Item(85, 11).Value = 730.2727
$ws.Cells.Item(85, 12).Value = 1194.75
$ws.Cells.Item(85, 13).Value = 517.7273
$ws.Cells.Item(85, 14).Value = -3690.75
$ws.Cells.Item(113, 8).Value = 64550.75
$ws.Cells.Item(113, 9).Value = 64550.75
$ws.Cells.Item(113, 11).Value = 64550.75
$ws.Cells.Item(113, 13).Value = -62380.75
$ws.Cells.Item(126, 8).Value = 3022
$ws.Cells.Item(126, 9).Value = 3329.3333
$ws.Cells.Item(126, 11).Value = 9987.999899999999
$ws.Cells.Item(126, 13).Value = -7517.999899999999
$ws.Cells.Item(132, 8).Value = 13220.107
$ws.Cells.Item(132, 9).Value = 12612.529
$ws.Cells.Item(132, 11).Value = 37837.587
$ws.Cells.Item(132, 13).Value = -35307.587
$ws.Cells.Item(136, 8).Value = 4190.25
$ws.Cells.Item(136, 9).Value = 5981.6665
$ws.Cells.Item(136, 11).Value = 17944.9995
$ws.Cells.Item(136, 13).Value = -15394.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 19998.2
$ws.Cells.Item(14, 10).Value = 19998.2
$ws.Cells.Item(14, 12).Value = 19998.2
$ws.Cells.Item(14, 14).Value = -20334.2
$ws.Cells.Item(136, 8).Value = 4510.086
$ws.Cells.Item(136, 9).Value = 2197.52
$ws.Cells.Item(136, 10).Value = 10291.5
$ws.Cells.Item(136, 11).Value = 6592.559999999999
$ws.Cells.Item(136, 12).Value = 30874.5
$ws.Cells.Item(136, 13).Value = -4042.559999999999
$ws.Cells.Item(136, 14).Value = -35974.5

